$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The columns BP:BU keep their position, but the "Usage" series that lives in
# each column is re-shuffled: both the header label (row 1) and the data
# value (row 2) for a given series move together to their new column.
#
#   column  old header       old value  ->  new header       new value
#   BP      tkm-N2Usage      26.2       ->  tkm-N3Usage      123.8
#   BQ      tkm-N3Usage      123.8      ->  tkm-N1Usage      8
#   BR      tkm-SZMUsage     388.5      ->  pkmUsage         838.8
#   BS      pkmUsage         838.8      ->  tkm-N2Usage      26.2
#   BT      tkm-N1Usage      8          ->  keroseneUsage    121
#   BU      keroseneUsage    121        ->  tkm-SZMUsage     388.5

$ws.Range("BP1").Value = "tkm-N3Usage"
$ws.Range("BQ1").Value = "tkm-N1Usage"
$ws.Range("BR1").Value = "pkmUsage"
$ws.Range("BS1").Value = "tkm-N2Usage"
$ws.Range("BT1").Value = "keroseneUsage"
$ws.Range("BU1").Value = "tkm-SZMUsage"

$ws.Range("BP2").Value = 123.8
$ws.Range("BQ2").Value = 8
$ws.Range("BR2").Value = 838.8
$ws.Range("BS2").Value = 26.2
$ws.Range("BT2").Value = 121
$ws.Range("BU2").Value = 388.5
